$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The stats rows for this player were re-ordered / updated as part of
# "updated activity till excel form". Columns are:
#   C = runs, D = balls, E = fours, F = sixes
# Keep values stored as text (matches the original numberStoredAsText layout).
$ws.Range("C2:F5").NumberFormat = "@"

$ws.Range("C2").Value = "15"
$ws.Range("D2").Value = "16"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "0"

$ws.Range("C3").Value = "10"
$ws.Range("D3").Value = "9"

$ws.Range("C4").Value = "4"
$ws.Range("D4").Value = "1"

$ws.Range("C5").Value = "54"
$ws.Range("D5").Value = "36"
$ws.Range("E5").Value = "2"
$ws.Range("F5").Value = "3"
